$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.012.92"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.654.68"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.85"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3908"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3850"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.27"
$ws.Range("E9").Value = "  +3.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.359"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08454"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.07"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.135"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.887"
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001318"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "1.653.16"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.73"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06986"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.84"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.929"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "23.978.47"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.029"
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.14"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.49"
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.449"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.55"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.909"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.489"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "1.836.92"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.046"
$ws.Range("E34").Value = "  +6.88%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02979"
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.780"
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.90"
$ws.Range("E38").Value = "  +5.59%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09158"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7571"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.48"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.428"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.33"
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6966"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.459"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.085"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9992"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08303"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.77"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.217"
$ws.Range("E51").Value = "  +1.44%  "
